$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.150.98'
$ws.Range("E2").Value = '  +1.39%  '

$ws.Range("D3").Value = '1.907.40'
$ws.Range("E3").Value = '  +1.73%  '

$ws.Range("D4").Value = '''1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").Value = '''327.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.69%  '

$ws.Range("D6").Value = '''1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Value = '''0.4603'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").Value = '''0.3934'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.72%  '

$ws.Range("D9").Value = '''46.71'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.48%  '

$ws.Range("D10").Value = '''0.07933'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.98%  '

$ws.Range("D11").Value = '''1.002'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.25%  '

$ws.Range("E12").Value = '  +2.06%  '

$ws.Range("D13").Value = '1.915.07'
$ws.Range("E13").Value = '  +2.35%  '

$ws.Range("D14").Value = '''7.092'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.39%  '

$ws.Range("E15").Value = '  +0.85%  '

$ws.Range("D16").Value = '''0.06951'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.09%  '

$ws.Range("D17").Value = '''88.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.10%  '

$ws.Range("D18").Value = '''1.004'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.08%  '

$ws.Range("E19").Value = '  +0.27%  '

$ws.Range("D20").Value = '''17.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.97%  '

$ws.Range("D22").Value = '29.155.65'
$ws.Range("E22").Value = '  +1.39%  '

$ws.Range("E23").Value = '  +1.50%  '

$ws.Range("D24").Value = '''11.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.44%  '

$ws.Range("D25").Value = '2.129.31'
$ws.Range("E25").Value = '  +1.53%  '

$ws.Range("E26").Value = '  -2.09%  '

$ws.Range("D27").Value = '''156.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.54%  '

$ws.Range("D28").Value = '''19.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.01%  '

$ws.Range("D29").Value = '''6.149'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.30%  '

$ws.Range("D30").Value = '''1.992'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.03%  '

$ws.Range("E31").Value = '  -0.25%  '

$ws.Range("D32").Value = '''0.09372'
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").Value = '''0.9248'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.56%  '

$ws.Range("D34").Value = '''5.331'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.53%  '

$ws.Range("E35").Value = '  +0.82%  '

$ws.Range("D36").Value = '''3.274'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.46%  '

$ws.Range("D37").Value = '''1.218'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.59%  '

$ws.Range("D38").Value = '''0.05833'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.12%  '

$ws.Range("D39").Value = '''0.02103'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.24%  '

$ws.Range("D40").Value = '''7.928'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.19%  '

$ws.Range("E41").Value = '  -0.25%  '

$ws.Range("D42").Value = '''0.5745'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.87%  '

$ws.Range("D43").Value = '''0.1799'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.67%  '

$ws.Range("D44").Value = '''9.944'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.68%  '

$ws.Range("D45").Value = '''2.244'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.89%  '

$ws.Range("D46").Value = '''11.94'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.65%  '

$ws.Range("D47").Value = '''0.5405'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.13%  '

$ws.Range("D48").Value = '''0.07076'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.92%  '

$ws.Range("D49").Value = '''1.875'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.56%  '

$ws.Range("D50").Value = '''2.552'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.85%  '

$ws.Range("D51").Value = '''112.86'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.58%  '
